# fix bug get list of campus
# Corrects the student roll-number/name and the "Content" (class/subject)
# values that were pulled from the wrong source rows, and fixes a couple
# of mismatched DOB/Sex/SubjectCode/AvgMark values on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll number + full name were wrong for every data row (B2:C6)
$ws.Range("B2:B6").Value = "HE130585"
$ws.Range("C2:C6").Value = "Nguyễn Ngọc Hải"

# Content column (K) per row
$ws.Range("K2").Value = "Đàn tranh"
$ws.Range("K3").Value = "Vovinam 1"
$ws.Range("K5").Value = "Vovinam 2"
$ws.Range("K6").Value = "Vovinam 3"

# Row 6 had several fields copied from the wrong record
$ws.Range("D6").Value = "23/07/1999"
$ws.Range("E6").Value = "Nữ"
$ws.Range("I6").Value = "ĐTR101"
$ws.Range("J6").Value = 9

# Column M width was adjusted while fixing the sheet
$ws.Range("M1").ColumnWidth = 20.6

# Leave the cursor where the author last left it when saving
[void]$ws.Range("L11").Select()
